# Update US model copy
# - Set the "motorbikes" row (row 6) to be subject to LCFS for both
#   passenger (B) and freight (C) vehicle types on the PVTStL sheet.
# - Move the saved cell selection on that sheet from B5 to D6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PVTStL")

# Flip the "motorbikes" row from 0/0 to 1/1.
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 1

# Update the sheet's remembered selection to D6 without stealing the
# active-tab state from whichever sheet is currently active.
$previouslyActive = $wb.ActiveSheet
$ws.Range("D6").Select()
$previouslyActive.Activate()
